$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "73.555.30"
$ws.Range("E2").Value = "  +2.00%  "
$ws.Range("D3").Value = "4.000.14"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'614.07"
$ws.Range("E5").Value = "  +13.89%  "
$ws.Range("D6").Value = "'166.20"
$ws.Range("E6").Value = "  +11.37%  "
$ws.Range("D7").Value = "'0.684"
$ws.Range("E7").Value = "  -1.54%  "
$ws.Range("D8").Value = "'0.998"
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("D10").Value = "'0.190"
$ws.Range("E10").Value = "  +10.45%  "
$ws.Range("D11").Value = "'56.36"
$ws.Range("E11").Value = "  +6.06%  "
$ws.Range("D12").Value = "'0.0000343"
$ws.Range("E12").Value = "  +3.86%  "
$ws.Range("D13").Value = "'11.11"
$ws.Range("E13").Value = "  +2.18%  "
$ws.Range("D14").Value = "4.624.00"
$ws.Range("D15").Value = "3.972.94"
$ws.Range("E15").Value = "  -1.96%  "
$ws.Range("E16").Value = "  +2.91%  "
$ws.Range("D17").Value = "'14.18"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "'20.56"
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("D19").Value = "73.252.92"
$ws.Range("E19").Value = "  +1.71%  "
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").Value = "'442.09"
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("D22").Value = "'4.88"
$ws.Range("E22").Value = "  +14.09%  "
$ws.Range("E23").Value = "  -2.05%  "
$ws.Range("D24").Value = "'3.38"
$ws.Range("E24").Value = "  -3.22%  "
$ws.Range("D25").Value = "'14.15"
$ws.Range("E25").Value = "  -3.17%  "
$ws.Range("D26").Value = "'4.05"
$ws.Range("E26").Value = "  -7.33%  "
$ws.Range("D27").Value = "'11.11"
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("D28").Value = "'5.93"
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("D29").Value = "'10.49"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("D30").Value = "'36.23"
$ws.Range("E30").Value = "  -2.29%  "
$ws.Range("D31").Value = "'7.73"
$ws.Range("E31").Value = "  -7.75%  "
$ws.Range("D32").Value = "'13.70"
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("E33").Value = "  -2.99%  "
$ws.Range("D36").Value = "0.0₃0999"
$ws.Range("E36").Value = "  +9.48%  "
$ws.Range("D37").Value = "'639.51"
$ws.Range("E37").Value = "  -6.08%  "
$ws.Range("D38").Value = "'0.432"
$ws.Range("E38").Value = "  -5.25%  "
$ws.Range("D39").Value = "'3.40"
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "'11.13"
$ws.Range("E42").Value = "  -1.52%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "'3.29"
$ws.Range("E44").Value = "  -4.68%  "
$ws.Range("D45").Value = "'0.0485"
$ws.Range("E45").Value = "  -1.50%  "
$ws.Range("E46").Value = "  -1.11%  "
$ws.Range("D47").Value = "'3.40"
$ws.Range("E47").Value = "  +2.75%  "
$ws.Range("D48").Value = "'2.60"
$ws.Range("E48").Value = "  -1.49%  "
$ws.Range("D49").Value = "'2.87"
$ws.Range("E49").Value = "  +27.83%  "
$ws.Range("D50").Value = "2.849.45"
$ws.Range("E50").Value = "  +2.43%  "

# Row 34/35 identity swap (InjectiveProtocol <-> OKB)
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "'71.37"
$ws.Range("E34").Value = "  +6.82%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "'47.76"
$ws.Range("E35").Value = "  -3.69%  "

# Row 51 replacement (Stacks -> FLOKI)
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").Value = "'0.000280"
$ws.Range("E51").Value = "  +2.32%  "
